# Add newly-scraped PMC article URLs to the "urls" sheet.
# Rows 1-6 (header + first 5 URLs) are already present and untouched.
# We append rows 7-19 with the new URLs, re-creating hyperlinks (and the
# "Hyperlink" cell style) only where the original commit added one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: new URL, styled like a hyperlink cell but with no live link.
$ws.Range("A7").Value = "https://www.ncbi.nlm.nih.gov/pmc/articles/PMC10977893/"
$ws.Range("A7").Style = "Hyperlink"

# Row 8: new URL with a real hyperlink.
$ws.Range("A8").Value = "https://www.ncbi.nlm.nih.gov/pmc/articles/PMC8446952/"
$ws.Hyperlinks.Add($ws.Range("A8"), "https://www.ncbi.nlm.nih.gov/pmc/articles/PMC8446952/") | Out-Null
$ws.Range("A8").Style = "Hyperlink"

# Row 9: duplicate of the row-5 URL, plain (no style, no hyperlink).
$ws.Range("A9").Value = "https://www.ncbi.nlm.nih.gov/pmc/articles/PMC9971900/"

# Row 10: new URL with a real hyperlink.
$ws.Range("A10").Value = "https://www.ncbi.nlm.nih.gov/pmc/articles/PMC11336597/"
$ws.Hyperlinks.Add($ws.Range("A10"), "https://www.ncbi.nlm.nih.gov/pmc/articles/PMC11336597/") | Out-Null
$ws.Range("A10").Style = "Hyperlink"

# Row 11: new URL with a real hyperlink.
$ws.Range("A11").Value = "https://www.ncbi.nlm.nih.gov/pmc/articles/PMC10558031/"
$ws.Hyperlinks.Add($ws.Range("A11"), "https://www.ncbi.nlm.nih.gov/pmc/articles/PMC10558031/") | Out-Null
$ws.Range("A11").Style = "Hyperlink"

# Row 12-17: plain new URLs, no style, no hyperlink.
$ws.Range("A12").Value = "https://www.ncbi.nlm.nih.gov/pmc/articles/PMC6655584/"
$ws.Range("A13").Value = "https://www.ncbi.nlm.nih.gov/pmc/articles/PMC10118127/"
$ws.Range("A14").Value = "https://www.ncbi.nlm.nih.gov/pmc/articles/PMC11133169/"
$ws.Range("A15").Value = "https://www.ncbi.nlm.nih.gov/pmc/articles/PMC10815757/"
$ws.Range("A16").Value = "https://www.ncbi.nlm.nih.gov/pmc/articles/PMC9761729/"
$ws.Range("A17").Value = "https://www.ncbi.nlm.nih.gov/pmc/articles/PMC6156563/"

# Row 18: new URL with a real hyperlink.
$ws.Range("A18").Value = "https://www.ncbi.nlm.nih.gov/pmc/articles/PMC9324041/"
$ws.Hyperlinks.Add($ws.Range("A18"), "https://www.ncbi.nlm.nih.gov/pmc/articles/PMC9324041/") | Out-Null
$ws.Range("A18").Style = "Hyperlink"

# Row 19: plain new URL, no style, no hyperlink.
$ws.Range("A19").Value = "https://www.ncbi.nlm.nih.gov/pmc/articles/PMC10463182/"

# Keep the active-cell selection in sync with where Excel would have left it.
$ws.Range("D21").Select() | Out-Null
